# Applies the "Updated cryptos list" price/volume refresh to Sheet1.
# Source data model: column D (Price) and column E (Volume(1h)) are
# plain text cells (inlineStr) in the workbook even when the price text
# looks numeric (e.g. "567.98"), because thousands separators use "."
# and some values carry significant trailing zeros ("27.90", "8.00").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new Price values (column D) parse as plain numbers
# (e.g. 567.60, 27.90). Writing them straight through Range.Value would
# let Excel auto-convert the cell to a Number and silently drop the
# trailing zero (567.60 -> 567.6), which would not match the source data.
# Force those specific cells to Text format first so the assignment below
# keeps them as literal strings, exactly like the rest of the column.
foreach ($addr in @("D5", "D6", "D9", "D12", "D13", "D18", "D19", "D20", "D21", "D23", "D24", "D26", "D27", "D28", "D29", "D32", "D35", "D36", "D37", "D38", "D39", "D40", "D43", "D44", "D45", "D47", "D49", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row-by-row Price (D) / Volume(1h) (E) updates
$ws.Range("D2").Value = "62.690.17"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "2.577.45"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("D5").Value = "567.60"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Value = "153.07"
$ws.Range("E6").Value = "  -2.64%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -2.50%  "
$ws.Range("D9").Value = "0.114"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "0.375"
$ws.Range("E12").Value = "  -2.82%  "
$ws.Range("D13").Value = "27.90"
$ws.Range("E13").Value = "  -2.02%  "
$ws.Range("D14").Value = "3.045.10"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("D16").Value = "62.636.87"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "2.610.84"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "11.86"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "7.42"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("D20").Value = "4.41"
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").Value = "335.85"
$ws.Range("E21").Value = "  -2.65%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "67.04"
$ws.Range("E23").Value = "  -0.83%  "
$ws.Range("D24").Value = "1.85"
$ws.Range("E24").Value = "  +5.22%  "
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("D26").Value = "1.61"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").Value = "9.03"
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").Value = "562.86"
$ws.Range("E28").Value = "  -6.28%  "
$ws.Range("D29").Value = "8.02"
$ws.Range("E29").Value = "  +1.15%  "
$ws.Range("E30").Value = "  +0.31%  "
$ws.Range("E31").Value = "  -3.14%  "
$ws.Range("D32").Value = "1.99"
$ws.Range("E32").Value = "  -4.37%  "
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("E34").Value = "  -3.53%  "
$ws.Range("D35").Value = "5.26"
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "0.395"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("D38").Value = "19.34"
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("D39").Value = "153.84"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "1.85"
$ws.Range("E40").Value = "  -2.22%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +2.05%  "
$ws.Range("D43").Value = "157.52"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("D44").Value = "23.47"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "3.83"
$ws.Range("E45").Value = "  -3.04%  "
$ws.Range("E46").Value = "  -3.96%  "
$ws.Range("D47").Value = "0.618"
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("E48").Value = "  -3.64%  "
$ws.Range("D49").Value = "0.0241"
$ws.Range("E49").Value = "  -3.86%  "
$ws.Range("E50").Value = "  -2.54%  "
$ws.Range("D51").Value = "0.767"
$ws.Range("E51").Value = "  -1.64%  "
